$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Move to location (11, 8) and remove the toolkit."
$ws.Range("B1").Value = "['Robot32']"
$ws.Range("E1").Value = "(11, 8)"

# Row 2
$ws.Range("A2").Value = "Move to location (7, 5) and remove the liquid spill."
$ws.Range("E2").Value = "(7, 5)"

# Row 3
$ws.Range("A3").Value = "Move to location (8, 6) and remove the large debris."
$ws.Range("B3").Value = "['Robot42', 'Robot29']"
$ws.Range("E3").Value = "(8, 6)"

# Row 4
$ws.Range("A4").Value = "Move to location (2, 4) and remove the dust."
$ws.Range("B4").Value = "['Robot50', 'Robot28']"
$ws.Range("E4").Value = "(2, 4)"

# Row 5
$ws.Range("A5").Value = "Move to location (5, 2) and remove the grass."
$ws.Range("B5").Value = "['Robot41']"
$ws.Range("E5").Value = "(5, 2)"

# Row 6
$ws.Range("A6").Value = "Move to location (6, 7) and remove the small debris."
$ws.Range("E6").Value = "(6, 7)"

# Row 7
$ws.Range("A7").Value = "Move to location (3, 6) and remove the vehicle."
$ws.Range("B7").Value = "['Robot13']"
$ws.Range("E7").Value = "(3, 6)"

# Row 8
$ws.Range("A8").Value = "Move to location (6, 6) and remove the construction materials."
$ws.Range("B8").Value = "['Robot2', 'Robot9', 'Robot23']"
$ws.Range("E8").Value = "(6, 6)"

# Row 9
$ws.Range("A9").Value = "Move to location (3, 9) and remove the tree branches."
$ws.Range("E9").Value = "(3, 9)"

# Row 10
$ws.Range("A10").Value = "Move to location (6, 6) and remove the screws."
$ws.Range("B10").Value = "['Robot15']"
$ws.Range("E10").Value = "(6, 6)"
